$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rodada 1")

# Row 2
$ws.Range("D2").Value = 49
$ws.Range("F2").Value = 47.86

# Row 3
$ws.Range("D3").Value = 54.6

# Row 4
$ws.Range("D4").Value = 84.86
$ws.Range("F4").Value = 49.76

# Row 5
$ws.Range("D5").Value = 47.86
$ws.Range("F5").Value = 72.7

# Row 6
$ws.Range("D6").Value = 60.16
$ws.Range("F6").Value = 72.86

# Row 7
$ws.Range("D7").Value = 74.06
$ws.Range("F7").Value = 61.56

# Row 8
$ws.Range("D8").Value = 73.95999999999999
$ws.Range("F8").Value = 53.66

# Row 9
$ws.Range("D9").Value = 64.7
$ws.Range("F9").Value = 68.06

# Row 10
$ws.Range("D10").Value = 63.9
$ws.Range("F10").Value = 84.26000000000001

# Row 11
$ws.Range("D11").Value = 54.16
$ws.Range("F11").Value = 72.45

# Row 12
$ws.Range("D12").Value = 38.26
$ws.Range("F12").Value = 62.56

# Row 13
$ws.Range("D13").Value = 60.2
$ws.Range("F13").Value = 59.25

# Row 14
$ws.Range("D14").Value = 54.1
$ws.Range("F14").Value = 61.96

# Row 15
$ws.Range("D15").Value = 73.76000000000001
$ws.Range("F15").Value = 57.6

# Row 16
$ws.Range("D16").Value = 42.96
$ws.Range("F16").Value = 57.45

# Row 17
$ws.Range("D17").Value = 68.06
$ws.Range("F17").Value = 58.96
